$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Deadline column (D) for rows 2-7 gets a date-like value that was
# typed with a leading apostrophe (quote prefix) so Excel stores it as
# literal text "20/09/2024'" even though the column's number format is
# a date format. Assigning via .Value keeps the cells on their existing
# (date-formatted) style instead of Excel deriving a brand-new style.
$ws.Range("D2").Value = "20/09/2024'"
$ws.Range("D3").Value = "20/09/2024'"
$ws.Range("D4").Value = "20/09/2024'"
$ws.Range("D5").Value = "20/09/2024'"
$ws.Range("D6").Value = "20/09/2024'"
$ws.Range("D7").Value = "20/09/2024'"

# Update the view: scroll back to A1 (no frozen/shifted topLeftCell) and
# move the selection to G19.
$ws.Range("A1").Select() | Out-Null
$ws.Range("G19").Select() | Out-Null
